$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Percent Error"
$ws.Range("C9").Value = 0.02
$ws.Range("C9").NumberFormat = "0%"

$ws.Range("C18").Select()
